# Auto-generated edit script: updates market-price / profit columns (H-N)
# on the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to match the scheduled-runner
# data refresh described in the commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 389.6316
$ws.Range("I53").Value = 340.16666
$ws.Range("J53").Value = 474.42856
$ws.Range("K53").Value = 340.16666
$ws.Range("L53").Value = 474.42856
$ws.Range("M53").Value = 296.83334
$ws.Range("N53").Value = -1748.42856
$ws.Range("H64").Value = 103000
$ws.Range("I64").Value = 169166.67
$ws.Range("J64").Value = 3750
$ws.Range("K64").Value = 169166.67
$ws.Range("L64").Value = 3750
$ws.Range("M64").Value = -168918.67
$ws.Range("N64").Value = -4246
$ws.Range("H67").Value = 103000
$ws.Range("I67").Value = 169166.67
$ws.Range("J67").Value = 3750
$ws.Range("K67").Value = 169166.67
$ws.Range("L67").Value = 3750
$ws.Range("M67").Value = -168308.67
$ws.Range("N67").Value = -5466
$ws.Range("H132").Value = 6585050
$ws.Range("I132").Value = 6950608.5
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 20851825.5
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -20849295.5
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6324.554
$ws.Range("I32").Value = 5444.3545
$ws.Range("J32").Value = 23708.5
$ws.Range("K32").Value = 5444.3545
$ws.Range("L32").Value = 23708.5
$ws.Range("M32").Value = -5157.3545
$ws.Range("N32").Value = -24282.5
$ws.Range("H45").Value = 64032.25
$ws.Range("I45").Value = 101110.7
$ws.Range("J45").Value = 2234.8333
$ws.Range("K45").Value = 101110.7
$ws.Range("L45").Value = 2234.8333
$ws.Range("M45").Value = -100733.7
$ws.Range("N45").Value = -2988.8333
$ws.Range("H80").Value = 27441.2
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 27441.2
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 27441.2
$ws.Range("M80").Value = ""
$ws.Range("N80").Value = -29437.2
$ws.Range("H82").Value = 33987.25
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 33987.25
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 33987.25
$ws.Range("N82").Value = -34709.25
$ws.Range("H83").Value = 27441.2
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 27441.2
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 82323.60000000001
$ws.Range("M83").Value = ""
$ws.Range("N83").Value = -92307.60000000001
$ws.Range("H85").Value = 33987.25
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 33987.25
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 33987.25
$ws.Range("N85").Value = -36483.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 37666.645
$ws.Range("I20").Value = 47564.41
$ws.Range("J20").Value = 1374.8334
$ws.Range("K20").Value = 47564.41
$ws.Range("L20").Value = 1374.8334
$ws.Range("M20").Value = -47317.41
$ws.Range("N20").Value = -1868.8334
$ws.Range("H80").Value = 2187.4614
$ws.Range("I80").Value = 750
$ws.Range("J80").Value = 2717.0527
$ws.Range("K80").Value = 750
$ws.Range("L80").Value = 2717.0527
$ws.Range("M80").Value = 248
$ws.Range("N80").Value = -4713.0527
$ws.Range("H83").Value = 2187.4614
$ws.Range("I83").Value = 750
$ws.Range("J83").Value = 2717.0527
$ws.Range("K83").Value = 3750
$ws.Range("L83").Value = 13585.2635
$ws.Range("M83").Value = 1242
$ws.Range("N83").Value = -23569.2635
$ws.Range("H99").Value = 1328.326
$ws.Range("I99").Value = 908.2857
$ws.Range("J99").Value = 1681.16
$ws.Range("K99").Value = 908.2857
$ws.Range("L99").Value = 1681.16
$ws.Range("M99").Value = 589.7143
$ws.Range("N99").Value = -4677.16
$ws.Range("H105").Value = 101834
$ws.Range("I105").Value = 201754.8
$ws.Range("J105").Value = 1913.2
$ws.Range("K105").Value = 201754.8
$ws.Range("L105").Value = 1913.2
$ws.Range("M105").Value = -200007.8
$ws.Range("N105").Value = -5407.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1333
$ws.Range("I16").Value = 1333
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1333
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1046
$ws.Range("N16").Value = ""
$ws.Range("H113").Value = 1333
$ws.Range("I113").Value = 1333
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1333
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 837
$ws.Range("N113").Value = ""
$ws.Range("H132").Value = 3985.577
$ws.Range("I132").Value = 3195.4211
$ws.Range("J132").Value = 6130.2856
$ws.Range("K132").Value = 9586.263300000001
$ws.Range("L132").Value = 18390.8568
$ws.Range("M132").Value = -7056.263300000001
$ws.Range("N132").Value = -23450.8568
$ws.Range("H134").Value = 1538
$ws.Range("I134").Value = 1305.3334
$ws.Range("J134").Value = 3166.6667
$ws.Range("K134").Value = 3916.0002
$ws.Range("L134").Value = 9500.000100000001
$ws.Range("M134").Value = -1381.0002
$ws.Range("N134").Value = -14570.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1652.8
$ws.Range("I97").Value = 1388.8
$ws.Range("J97").Value = 1916.8
$ws.Range("K97").Value = 4166.4
$ws.Range("L97").Value = 5750.4
$ws.Range("M97").Value = -3670.4
$ws.Range("N97").Value = -6742.4
$ws.Range("H107").Value = 605019.0600000001
$ws.Range("I107").Value = 960
$ws.Range("J107").Value = 772813.25
$ws.Range("K107").Value = 2880
$ws.Range("L107").Value = 2318439.75
$ws.Range("M107").Value = -960
$ws.Range("N107").Value = -2322279.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H106").Value = 21000
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 21000
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 21000
$ws.Range("N106").Value = -23524

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 85583.25
$ws.Range("I40").Value = 501000
$ws.Range("J40").Value = 2499.9
$ws.Range("K40").Value = 501000
$ws.Range("L40").Value = 2499.9
$ws.Range("M40").Value = -500864
$ws.Range("N40").Value = -2771.9
$ws.Range("H105").Value = 44051.25
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 44051.25
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 44051.25
$ws.Range("N105").Value = -51039.25
$ws.Range("H122").Value = 2948.2222
$ws.Range("I122").Value = 2900.125
$ws.Range("J122").Value = 3333
$ws.Range("K122").Value = 8700.375
$ws.Range("L122").Value = 9999
$ws.Range("M122").Value = -6250.375
$ws.Range("N122").Value = -14899
$ws.Range("H136").Value = 2062.8333
$ws.Range("I136").Value = 2128.3333
$ws.Range("J136").Value = 1866.3334
$ws.Range("K136").Value = 6384.999899999999
$ws.Range("L136").Value = 5599.0002
$ws.Range("M136").Value = -3834.999899999999
$ws.Range("N136").Value = -10699.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6412706
$ws.Range("I62").Value = 76923070
$ws.Range("J62").Value = 2672.7273
$ws.Range("K62").Value = 76923070
$ws.Range("L62").Value = 2672.7273
$ws.Range("M62").Value = -76922446
$ws.Range("N62").Value = -3920.7273
$ws.Range("H65").Value = 6412706
$ws.Range("I65").Value = 76923070
$ws.Range("J65").Value = 2672.7273
$ws.Range("K65").Value = 384615350
$ws.Range("L65").Value = 13363.6365
$ws.Range("M65").Value = -384612230
$ws.Range("N65").Value = -19603.6365
$ws.Range("H100").Value = 72164.36
$ws.Range("I100").Value = 91409.17999999999
$ws.Range("J100").Value = 1600
$ws.Range("K100").Value = 182818.36
$ws.Range("L100").Value = 3200
$ws.Range("M100").Value = -182277.36
$ws.Range("N100").Value = -4282
$ws.Range("H107").Value = 250899.75
$ws.Range("I107").Value = 1199.6666
$ws.Range("J107").Value = 1000000
$ws.Range("K107").Value = 3598.9998
$ws.Range("L107").Value = 3000000
$ws.Range("M107").Value = -1678.9998
$ws.Range("N107").Value = -3003840
